# Weekly update: insert 3 new daily-price rows for Chirimoya (Mercado Mayorista
# Lo Valledor de Santiago) at the top of the data block (row 23), shifting all
# existing data rows down by 3. Excel's Insert() on the entire-row range takes
# care of shifting rows 23-109 down to 26-112 and growing the sheet dimension.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows above the current row 23 (pushes old rows 23..109 -> 26..112)
$ws.Range("A23:A25").EntireRow.Insert()

# --- New row 23 ---
$ws.Range("A23").Value = 6
$ws.Range("B23").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C23").Value = "Metropolitana"
$ws.Range("D23").Value = 44469
$ws.Range("E23").Value = 13
$ws.Range("F23").Value = "Fruta"
$ws.Range("G23").Value = 100107
$ws.Range("H23").Value = "Otros"
$ws.Range("I23").Value = 100107002
$ws.Range("J23").Value = "Chirimoya"
$ws.Range("K23").Value = "Cultivar IV Región"
$ws.Range("L23").Value = "Especial"
$ws.Range("M23").Value = 200
$ws.Range("N23").Value = 3100
$ws.Range("O23").Value = 3100
$ws.Range("P23").Value = 3100
$ws.Range("Q23").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R23").Value = "Provincia de Limarí"
$ws.Range("S23").Value = 3100
$ws.Range("T23").Value = 1

# --- New row 24 ---
$ws.Range("A24").Value = 6
$ws.Range("B24").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C24").Value = "Metropolitana"
$ws.Range("D24").Value = 44469
$ws.Range("E24").Value = 13
$ws.Range("F24").Value = "Fruta"
$ws.Range("G24").Value = 100107
$ws.Range("H24").Value = "Otros"
$ws.Range("I24").Value = 100107002
$ws.Range("J24").Value = "Chirimoya"
$ws.Range("K24").Value = "Cultivar IV Región"
$ws.Range("L24").Value = "Extra (doble especial)"
$ws.Range("M24").Value = 150
$ws.Range("N24").Value = 3400
$ws.Range("O24").Value = 3400
$ws.Range("P24").Value = 3400
$ws.Range("Q24").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R24").Value = "Provincia de Limarí"
$ws.Range("S24").Value = 3400
$ws.Range("T24").Value = 1

# --- New row 25 ---
$ws.Range("A25").Value = 6
$ws.Range("B25").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C25").Value = "Metropolitana"
$ws.Range("D25").Value = 44469
$ws.Range("E25").Value = 13
$ws.Range("F25").Value = "Fruta"
$ws.Range("G25").Value = 100107
$ws.Range("H25").Value = "Otros"
$ws.Range("I25").Value = 100107002
$ws.Range("J25").Value = "Chirimoya"
$ws.Range("K25").Value = "Cultivar IV Región"
$ws.Range("L25").Value = "Primera"
$ws.Range("M25").Value = 200
$ws.Range("N25").Value = 2800
$ws.Range("O25").Value = 2800
$ws.Range("P25").Value = 2800
$ws.Range("Q25").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R25").Value = "Provincia de Limarí"
$ws.Range("S25").Value = 2800
$ws.Range("T25").Value = 1
